# 15/12/2023 * Updates to characterise the delay correction between the
# Synchronisation and the pixels sequence
#
# Appends 9 new test-result rows (39-47) to the "no_regression" sheet,
# covering DRE_DMX_UT_5031a-g and DRE_DMX_UT_5032a-b.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- DRE_DMX_UT_5031a..e (rows 39-43) ------------------------------------
$ws.Range("A39").Value = "DRE_DMX_UT_5031a"
$ws.Range("A40").Value = "DRE_DMX_UT_5031b"
$ws.Range("A41").Value = "DRE_DMX_UT_5031c"
$ws.Range("A42").Value = "DRE_DMX_UT_5031d"
$ws.Range("A43").Value = "DRE_DMX_UT_5031e"

$ws.Range("B39").Value = "Simulation with g_RA_DELAY=0"
$ws.Range("B40").Value = "Simulation with g_RA_DELAY=1"
$ws.Range("B41").Value = "Simulation with g_RA_DELAY=2"
$ws.Range("B42").Value = "Simulation with g_RA_DELAY=3"
$ws.Range("B43").Value = "Simulation with g_RA_DELAY=4"

$ws.Range("C39").Value = "PASS"
$ws.Range("C40").Value = "PASS"
$ws.Range("C41").Value = "PASS"
$ws.Range("C42").Value = "PASS"
$ws.Range("C43").Value = "PASS"

$ws.Range("D39").Value = "RAS / DEMUX delay is not correct"
$ws.Range("D40").Value = "RAS / DEMUX delay is not correct"
$ws.Range("D41").Value = "RAS / DEMUX delay is not correct"
$ws.Range("D42").Value = "RAS / DEMUX delay is not correct"
$ws.Range("D43").Value = "RAS / DEMUX delay is not correct"

# --- DRE_DMX_UT_5031f..g (rows 44-45) ------------------------------------
$ws.Range("A44").Value = "DRE_DMX_UT_5031f"
$ws.Range("A45").Value = "DRE_DMX_UT_5031g"

$ws.Range("B44").Value = "Simulation with g_RA_DELAY=5"
$ws.Range("B45").Value = "Simulation with g_RA_DELAY=6"

$ws.Range("C44").Value = "PASS"
$ws.Range("C45").Value = "PASS"

$ws.Range("D44").Value = "RAS / DEMUX delay is not correct"
$ws.Range("D45").Value = "RAS / DEMUX delay is not correct"

# --- DRE_DMX_UT_5032a..b (rows 46-47) ------------------------------------
$ws.Range("A46").Value = "DRE_DMX_UT_5032a"
$ws.Range("A47").Value = "DRE_DMX_UT_5032b"

$ws.Range("C46").Value = "PASS"
$ws.Range("C47").Value = "PASS"

$ws.Range("D46").Value = "RAS / DEMUX delay is not correct"
$ws.Range("D47").Value = "RAS / DEMUX delay is correct"

$ws.Range("B46").Value = "Simulation with g_RA_DELAY=6, g_ERROR_DELAY=0"
$ws.Range("B47").Value = "Simulation with g_RA_DELAY=6, g_ERROR_DELAY=2"

# Column B holds long descriptive text; keep it wrapping like the rest of
# that column (style already used by rows 34-38).
$ws.Range("B39:B47").WrapText = $true

# --- View state: scroll down and select B50, as in the authored commit --
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 31
$ws.Range("B50").Select()
